# ---------------------------------------------------------------------------
# Add "2022-Q4" data to the workbook:
#   1. Insert a new worksheet named "2022-Q4" right after "总计", containing
#      the per-fund holdings table for that quarter (copy format from the
#      existing "2022-Q3" sheet, since the table layout/headers are the same).
#   2. Insert a new row at the top of the "总计" (summary) sheet's data with
#      the 2022-Q4 totals, pushing every other quarter down by one row.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item(1)          # "总计"
$q3Sheet    = $wb.Worksheets.Item(2)          # "2022-Q3" (soon to become index 3)

# ---------------------------------------------------------------------------
# 1. Build the new "2022-Q4" worksheet by duplicating the "2022-Q3" sheet
#    (same columns/headers/styles), then resizing it from 7 to 10 data rows
#    and overwriting every value.
# ---------------------------------------------------------------------------
$q3Sheet.Copy($null, $totalSheet)
$q4Sheet = $wb.Worksheets.Item(2)
$q4Sheet.Name = "2022-Q4"

# The copied sheet has rows 2-8 (7 funds). We need rows 2-11 (10 funds), so
# insert 3 more rows before row 9, copying the formatting from row 8.
$q4Sheet.Rows.Item(9).Resize(3).Insert()
$q4Sheet.Range("A8:H8").Copy()
$q4Sheet.Range("A9:H11").PasteSpecial(-4122)

# Columns B-G hold text values (fund codes/names/ratios formatted as text in
# the source data) - force a text number format so values aren't coerced to
# numbers (this also preserves leading zeros in fund codes).
$q4Sheet.Range("B2:G11").NumberFormat = "@"

$q4Data = @(
    @("257010", "国联安小盘精选混合",               "8.86", "74.78", "5.84", "0.5174", 3),
    @("006138", "国联安价值优选股票",               "0.61", "94.70", "5.47", "0.0334", 6),
    @("011243", "万家惠裕回报6个月持有期混合A",       "1.28", "29.05", "1.25", "0.0160", 4),
    @("009658", "汇丰晋信中小盘低波动策略股票A",      "0.85", "92.42", "1.54", "0.0131", 1),
    @("007288", "合煦智远消费主题股票C",             "0.11", "83.65", "2.88", "0.0032", 9),
    @("004791", "富荣中证500指数增强C",             "0.09", "90.60", "2.10", "0.0019", 7),
    @("011244", "万家惠裕回报6个月持有期混合C",       "0.10", "29.05", "1.25", "0.0012", 4),
    @("007287", "合煦智远消费主题股票A",             "0.03", "83.65", "2.88", "0.0009", 9),
    @("009775", "汇丰晋信中小盘低波动策略股票C",      "0.04", "92.42", "1.54", "0.0006", 1),
    @("004790", "富荣中证500指数增强A",             "0.02", "90.60", "2.10", "0.0004", 7)
)

for ($i = 0; $i -lt $q4Data.Count; $i++) {
    $row = $i + 2
    $entry = $q4Data[$i]
    $q4Sheet.Cells.Item($row, 1).Value = $i
    $q4Sheet.Cells.Item($row, 2).Value = $entry[0]
    $q4Sheet.Cells.Item($row, 3).Value = $entry[1]
    $q4Sheet.Cells.Item($row, 4).Value = $entry[2]
    $q4Sheet.Cells.Item($row, 5).Value = $entry[3]
    $q4Sheet.Cells.Item($row, 6).Value = $entry[4]
    $q4Sheet.Cells.Item($row, 7).Value = $entry[5]
    $q4Sheet.Cells.Item($row, 8).Value = $entry[6]
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert a new row 2 for 2022-Q4 and shift
#    the rest down (copy formatting from row 3 so borders/styles line up).
# ---------------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("A3:D3").Copy()
$totalSheet.Range("A2:D2").PasteSpecial(-4122)

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q4"
$totalSheet.Cells.Item(2, 3).Value = 10
$totalSheet.Cells.Item(2, 4).Value = 0.59
